# "Mise à jour de l'application"
# Adds the results of "Match Amical 2" (columns AM:AP) for the players
# who took part in it, mirroring the structure already used for
# "Match Amical 1" (columns AI:AL). AM = minutes played, AN = T/R
# (titulaire/remplaçant), AO = but (goal), AP = passe decisive (assist).
# All the dependent totals (B, F, G, H, I, Y, Z, AE, AF, PW, ...) are
# driven by existing formulas already on the sheet, so simply writing
# the raw inputs below is enough for Excel to recompute them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$matchAmical2 = @(
    @{ Row = 2;  Minutes = 45; TR = "T"; But = $null; Passe = $null }
    @{ Row = 3;  Minutes = 45; TR = "R"; But = $null; Passe = $null }
    @{ Row = 5;  Minutes = 45; TR = "T"; But = $null; Passe = $null }
    @{ Row = 6;  Minutes = 45; TR = "R"; But = $null; Passe = $null }
    @{ Row = 9;  Minutes = 45; TR = "R"; But = $null; Passe = $null }
    @{ Row = 10; Minutes = 45; TR = "R"; But = $null; Passe = $null }
    @{ Row = 11; Minutes = 45; TR = "T"; But = $null; Passe = $null }
    @{ Row = 12; Minutes = 45; TR = "T"; But = $null; Passe = $null }
    @{ Row = 13; Minutes = 30; TR = "R"; But = $null; Passe = $null }
    @{ Row = 14; Minutes = 45; TR = "R"; But = $null; Passe = $null }
    @{ Row = 15; Minutes = 45; TR = "T"; But = $null; Passe = $null }
    @{ Row = 16; Minutes = 45; TR = "R"; But = $null; Passe = 1 }
    @{ Row = 17; Minutes = 45; TR = "R"; But = 1;     Passe = $null }
    @{ Row = 18; Minutes = 45; TR = "T"; But = $null; Passe = $null }
    @{ Row = 19; Minutes = 45; TR = "R"; But = $null; Passe = $null }
    @{ Row = 20; Minutes = 45; TR = "T"; But = $null; Passe = $null }
    @{ Row = 21; Minutes = 15; TR = "R"; But = $null; Passe = $null }
    @{ Row = 22; Minutes = 45; TR = "T"; But = $null; Passe = $null }
    @{ Row = 24; Minutes = 45; TR = "T"; But = $null; Passe = $null }
    @{ Row = 25; Minutes = 45; TR = "R"; But = $null; Passe = $null }
    @{ Row = 26; Minutes = 45; TR = "T"; But = $null; Passe = $null }
)

foreach ($entry in $matchAmical2) {
    $r = $entry.Row
    $ws.Range("AM$r").Value = $entry.Minutes
    $ws.Range("AN$r").Value = $entry.TR
    if ($entry.But -ne $null) {
        $ws.Range("AO$r").Value = $entry.But
    }
    if ($entry.Passe -ne $null) {
        $ws.Range("AP$r").Value = $entry.Passe
    }
}

# Restore the frozen pane at column A (xSplit=1) and move the visible
# selection to where the author ended up after entering the data.
$win = $ws.Application.ActiveWindow
$win.FreezePanes = $false
$ws.Range("B1").Select()
$win.FreezePanes = $true
$ws.Range("AR15").Select()
